$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 3.8883659508071853
$ws.Range("C2").Value = 5.7366023908921173
$ws.Range("D2").Value = 3.4270214280954638
$ws.Range("E2").Value = 3.6509631212893767

$ws.Range("B3").Value = 6.5318022016907191
$ws.Range("C3").Value = 10.048852349745278
$ws.Range("D3").Value = 8.0593215041399819
$ws.Range("E3").Value = 1.8114155050670901

$excel.Goto($ws.Range("B1:E3"))
